# Updates cryptos list values (price + 1h volume change) to match latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text cells (coin name / link / percentage strings) - safe to assign directly.
$textUpdates = @{
    'E2' = '  -2.00%  '
    'E3' = '  +0.01%  '
    'E4' = '  -0.06%  '
    'E5' = '  -2.43%  '
    'E6' = '  -4.43%  '
    'E7' = '  +0.01%  '
    'E8' = '  +0.09%  '
    'E9' = '  -1.54%  '
    'E10' = '  -3.29%  '
    'E11' = '  -2.64%  '
    'E12' = '  -2.99%  '
    'E13' = '  -4.12%  '
    'E14' = '  -3.17%  '
    'E15' = '  -0.01%  '
    'E16' = '  +1.22%  '
    'E17' = '  -1.99%  '
    'E18' = '  -0.21%  '
    'E19' = '  -3.62%  '
    'E20' = '  -2.22%  '
    'E21' = '  -3.26%  '
    'E22' = '  -1.83%  '
    'E23' = '  -0.33%  '
    'E24' = '  -1.89%  '
    'E25' = '  -3.05%  '
    'E26' = '  +0.02%  '
    'E27' = '  -1.83%  '
    'E28' = '  -5.33%  '
    'E29' = '  +1.67%  '
    'E30' = '  -4.25%  '
    'E31' = '  -0.02%  '
    'E32' = '  -0.72%  '
    'E33' = '  -4.80%  '
    'E34' = '  -5.24%  '
    'E35' = '  -2.56%  '
    'E36' = '  -0.45%  '
    'E37' = '  -3.95%  '
    'E38' = '  -8.44%  '
    'E39' = '  -1.89%  '
    'E40' = '  -6.72%  '
    'B41' = 'dogwifhat'
    'C41' = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
    'E41' = '  -10.01%  '
    'E42' = '  -0.45%  '
    'B43' = 'Maker'
    'C43' = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
    'E43' = '  +1.39%  '
    'E44' = '  -6.19%  '
    'E45' = '  +0.07%  '
    'E47' = '  -5.41%  '
    'E48' = '  -2.65%  '
    'E49' = '  -0.62%  '
    'E50' = '  -8.33%  '
    'E51' = '  -0.38%  '
}

foreach ($cell in $textUpdates.Keys) {
    $ws.Range($cell).Value = $textUpdates[$cell]
}

# Price cells look numeric (e.g. "588.89", "2.907.14", "0.0000244") so Excel would
# otherwise coerce/round them as numbers. Force text format, assign, then restore the
# default "Normal" style so no stray number formatting is left behind.
$priceUpdates = @{
    'D2' = '63.110.37'
    'D3' = '3.139.27'
    'D5' = '588.89'
    'D6' = '137.28'
    'D8' = '3.136.26'
    'D9' = '0.516'
    'D10' = '0.145'
    'D11' = '5.25'
    'D12' = '0.456'
    'D13' = '0.0000244'
    'D14' = '34.12'
    'D15' = '3.655.05'
    'D17' = '63.086.83'
    'D18' = '3.127.68'
    'D19' = '6.66'
    'D20' = '472.33'
    'D21' = '14.13'
    'D22' = '0.698'
    'D23' = '7.66'
    'D24' = '85.36'
    'D25' = '12.96'
    'D27' = '2.71'
    'D28' = '7.92'
    'D29' = '2.09'
    'D30' = '6.86'
    'D32' = '26.70'
    'D34' = '2.53'
    'D35' = '1.07'
    'D36' = '52.33'
    'D37' = '5.76'
    'D38' = '0.0₃0690'
    'D39' = '0.0386'
    'D40' = '416.20'
    'D41' = '2.70'
    'D42' = '8.20'
    'D43' = '2.907.14'
    'D44' = '0.112'
    'D45' = '0.261'
    'D47' = '2.11'
    'D48' = '25.39'
    'D49' = '0.113'
    'D51' = '120.70'
}

foreach ($cell in $priceUpdates.Keys) {
    $range = $ws.Range($cell)
    $range.NumberFormat = '@'
    $range.Value = $priceUpdates[$cell]
    $range.Style = 'Normal'
}
